$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le 4; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Width = 100
    }
}
